$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

$dataCols = @("B","C","D","E","F","G","H","I","J","K")

# Row 17: columns B..K currently hold blank inline strings; the new event
# normalizes them to the literal text "nan" (matching every other
# already-logged row on this card).
foreach ($col in $dataCols) {
    $ws.Range($col + "17").Value = "nan"
}

# Row 18: brand-new service-log entry appended for Card19.
# Column A repeats the card number "19" as text (same as every other row),
# so force text interpretation with a quote prefix instead of letting Excel
# coerce it to a number.
$ws.Range("A18").Value = "'19"

# Columns B..K stay blank text cells (same shape as a freshly logged row,
# mirrored from how B17:K17 looked before being filled in) - force an empty
# text value rather than leaving the cell truly empty.
foreach ($col in $dataCols) {
    $ws.Range($col + "18").Value = "'"
}

$ws.Range("L18").Value = "23\8\2025"
$ws.Range("M18").Value = "797 t"
$ws.Range("N18").Value = "تم تغيير الجريده رقم 1"
$ws.Range("O18").Value = "الخبير"
